$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the existing row 1. This pushes the
# current rows 1-3 (with their values + styles) down to rows 2-4,
# matching the diff (old row1 -> new row2, old row2 -> new row3,
# old row3 -> new row4).
$ws.Rows("1:1").Insert()

# New duty-roster columns G:J carry extra names/numbers for each of
# the three shifted rows. Columns H and J hold phone numbers and must
# be formatted as text first so leading zeros survive.
$ws.Range("H2:H4").NumberFormat = "@"
$ws.Range("J2:J4").NumberFormat = "@"

$ws.Range("G2").Value = "Shaon"
$ws.Range("H2").Value = "0164567564"
$ws.Range("I2").Value = "Selim"
$ws.Range("J2").Value = "0154676341"

$ws.Range("G3").Value = "Faruk"
$ws.Range("H3").Value = "0164576654"
$ws.Range("I3").Value = "Faruk"
$ws.Range("J3").Value = "0146674647"

$ws.Range("G4").Value = "Rofik"
$ws.Range("H4").Value = "0154764797"
$ws.Range("I4").Value = "Moin"
$ws.Range("J4").Value = "0157866344"

# New column widths for the newly used columns H and J.
$ws.Columns("H:H").ColumnWidth = 11
$ws.Columns("J:J").ColumnWidth = 10

# Restore the selection Excel saved the file with.
$ws.Range("G10").Select() | Out-Null
